# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# "bff14131-..." file has been handed back (in sync with en-US) for the
# zh-cn language, and that both zh-cn and de-de now have a recorded
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: roll the new status up into the per-language columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: bff14131 row (row 2) has been handed back
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("J2").Value = "bff14131-34e5-4a0d-8767-4c0b91f10dc2.eed5eab7c750aa714842042a6d49d6126f8f23d9.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-01 02:48:49"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c6175a715b97a2e24f5106f0fe0ebe446d85c1b/e2e/bff14131-34e5-4a0d-8767-4c0b91f10dc2.md", "", "", "bff14131-34e5-4a0d-8767-4c0b91f10dc2.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c6175a715b97a2e24f5106f0fe0ebe446d85c1b/e2e/bff14131-34e5-4a0d-8767-4c0b91f10dc2.md", "", "", "bff14131-34e5-4a0d-8767-4c0b91f10dc2.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/d74fdfbc-439e-49e1-a535-1791de276e2d.md", "", "", "d74fdfbc-439e-49e1-a535-1791de276e2d.md")

# ---------------------------------------------------------------------
# de-de sheet: bff14131 row (row 2) also got a handback target recorded
# (status for de-de stays "Ready for handoff")
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("J2").Value = "bff14131-34e5-4a0d-8767-4c0b91f10dc2.eed5eab7c750aa714842042a6d49d6126f8f23d9.de-de.xlf"
$dede.Range("K2").Value = "2016-09-01 02:48:55"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c6175a715b97a2e24f5106f0fe0ebe446d85c1b/e2e/bff14131-34e5-4a0d-8767-4c0b91f10dc2.md", "", "", "bff14131-34e5-4a0d-8767-4c0b91f10dc2.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c6175a715b97a2e24f5106f0fe0ebe446d85c1b/e2e/bff14131-34e5-4a0d-8767-4c0b91f10dc2.md", "", "", "bff14131-34e5-4a0d-8767-4c0b91f10dc2.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/d74fdfbc-439e-49e1-a535-1791de276e2d.md", "", "", "d74fdfbc-439e-49e1-a535-1791de276e2d.md")

# ---------------------------------------------------------------------
# Column widths: widen the columns that now hold longer text, mirroring
# the autofit Excel performs after the new content is written.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.142857142857142
$overview.Columns.Item(6).ColumnWidth = 29.142857142857142

$zhcn.Columns.Item(3).ColumnWidth = 29.142857142857142
$zhcn.Columns.Item(9).ColumnWidth = 39.142857142857146
$zhcn.Columns.Item(10).ColumnWidth = 39.142857142857146

$dede.Columns.Item(3).ColumnWidth = 29.142857142857142
$dede.Columns.Item(9).ColumnWidth = 39.142857142857146
$dede.Columns.Item(10).ColumnWidth = 39.142857142857146

Write-Host "Handback report generated"
